$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'71.546.40"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +2.63%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'3.658.02"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +8.19%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.01%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'588.18"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +0.99%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'180.30"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +0.14%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'3.650.04"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  +8.12%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.624"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  +4.87%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = "'  +0.00%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'0.203"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +1.71%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.611"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +3.60%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'49.83"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  +2.85%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'0.0000286"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  +0.22%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('B14').Value = "'BitcoinCash"
$ws.Range('B14').Style = 'Normal'
$ws.Range('C14').Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range('C14').Style = 'Normal'
$ws.Range('D14').Value = "'682.00"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -0.05%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('B15').Value = "'WrappedliquidstakedEther2.0"
$ws.Range('B15').Style = 'Normal'
$ws.Range('C15').Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range('C15').Style = 'Normal'
$ws.Range('D15').Value = "'4.237.83"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +7.96%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'9.02"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +4.51%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'3.697.15"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +9.37%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'71.632.67"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +2.76%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = "'  +1.79%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'18.24"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  +3.10%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = "'  +3.22%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'0.941"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  +3.10%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'6.17"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +15.27%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'17.86"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +2.94%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'103.31"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +1.44%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = "'  +2.99%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'2.84"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +5.07%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'10.22"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  +4.79%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'  +4.92%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'9.18"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  +4.94%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'7.45"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +7.54%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'4.23"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +9.96%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'580.91"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  +4.57%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'11.34"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +2.32%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = "'  +2.41%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'59.54"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  +2.82%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'3.751.47"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  +3.96%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'  +0.12%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'  +2.93%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'35.63"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +0.68%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'0.0₃0765"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  +4.90%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'3.46"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +4.22%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('B43').Value = "'Fetch.AI"
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').Value = "'2.80"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +1.60%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('B44').Value = "'VeChain"
$ws.Range('B44').Style = 'Normal'
$ws.Range('C44').Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range('C44').Style = 'Normal'
$ws.Range('D44').Value = "'0.0462"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +8.25%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').Value = "'  +2.63%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'3.36"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -0.12%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'2.81"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +4.94%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'  +3.75%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = "'  +4.32%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'0.998"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  -0.22%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'134.02"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +3.01%  "
$ws.Range('E51').Style = 'Normal'
